# feat: add 2022-Q4 data
#
# 1) Insert a new row into the "总计" (summary) sheet for the 2022-Q4
#    quarter, pushing the existing 2022-Q2 / 2021-Q4 rows down, and bump
#    the 2022-Q2 holding-value total (0.2 -> 0.26).
# 2) Add a brand new "2022-Q4" worksheet (placed right after "总计",
#    before "2022-Q2") containing the quarter's fund holding breakdown.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" summary sheet: insert a new row 3 for "2022-Q4" and update
#    values.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Push row 3 ("2021-Q4") down to make room; this also copies row 2's
# cell formatting into the freshly inserted row.
$summary.Rows.Item(3).Insert()

# Row 2 becomes the new "2022-Q4" entry.
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.26

# Row 3 keeps the data that used to live in row 2 ("2022-Q2").
$summary.Range("A3").Value = 1
$summary.Range("A3").HorizontalAlignment = -4108
$summary.Range("A3").VerticalAlignment = -4160
$summary.Range("A3").Font.Bold = $true
$summary.Range("A3").Borders.LineStyle = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 4
$summary.Range("D3").Value = 0.2

# Row 4 ("2021-Q4") keeps its original values; only its row index
# label needs to move from 1 -> 2.
$summary.Range("A4").Value = 2

# ---------------------------------------------------------------------
# 2) New "2022-Q4" worksheet with the fund-holding breakdown.
# ---------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("总计")
$q4 = $wb.Worksheets.Add([Type]::Missing, $afterSheet)
$q4.Name = "2022-Q4"

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q4.Cells.Item(1, 2 + $i)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$rows = @(
    @(0, "012868", "易方达标普信息科技指数（QDII-LOF）人民币 C", "5.09", "91.36", "1.73", "0.0881", 9),
    @(1, "161128", "易方达标普信息科技指数（QDII-LOF）人民币",   "5.09", "91.36", "1.73", "0.0881", 9),
    @(2, "003721", "易方达标普信息科技指数（QDII-LOF）美元A",    "4.93", "91.36", "1.73", "0.0853", 9),
    @(3, "012869", "易方达标普信息科技指数（QDII-LOF）美元 C",   "0.16", "91.36", "1.73", "0.0028", 9)
)

$r = 2
foreach ($row in $rows) {
    $q4.Cells.Item($r, 1).Value = $row[0]
    $q4.Cells.Item($r, 1).HorizontalAlignment = -4108
    $q4.Cells.Item($r, 1).VerticalAlignment = -4160
    $q4.Cells.Item($r, 1).Font.Bold = $true
    $q4.Cells.Item($r, 1).Borders.LineStyle = 1

    $codeCell = $q4.Cells.Item($r, 2)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $row[1]

    $q4.Cells.Item($r, 3).Value = $row[2]

    for ($c = 4; $c -le 7; $c++) {
        $cell = $q4.Cells.Item($r, $c)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$c - 1]
    }

    $q4.Cells.Item($r, 8).Value = $row[7]
    $r++
}
